$wb = $excel.ActiveWorkbook

# --- Sheet "Contact Info" ---
$wsContact = $wb.Worksheets.Item("Contact Info")
# Phone number for Andrew Michel becomes a text value (with spaces) instead of a bare number
$wsContact.Range("B2").Value = "661 809 6224"
# Minor column width tweaks
$wsContact.Columns.Item(2).ColumnWidth = 15.642857142857142
$wsContact.Columns.Item(3).ColumnWidth = 24.142857142857142
# This sheet is no longer the active tab, leaving the cursor elsewhere
$wsContact.Activate() | Out-Null
$wsContact.Range("C33:D35").Select() | Out-Null

# --- Sheet "Rough Guide" ---
$wsRough = $wb.Worksheets.Item("Rough Guide")
# Expand the description of the base PacObject class
$wsRough.Range("A2").Value = "One main class which everything will inherit from. It can be called something like 'PacObject'. This will be the object type of our 2 dimensional array. PacObject will have a 'move' method. It's current x,y location (might not need this). Also some kind of detectHit method. We also need some sort of logic over what would 'win' in a collision. e.g. PacMan eats PacPrize but PacEnemy eats PacMan. Plus account for PacWall."
$wsRough.Rows.Item(2).RowHeight = 105
# New class note
$wsRough.Range("A7").Value = "PacWall inherits PacObject"
$wsRough.Rows.Item(5).RowHeight = 45
$wsRough.Rows.Item(10).RowHeight = 60
# New note about a central clock for managing movements
$wsRough.Range("A14").Value = "Since we're doing array based collisions we need to determine some kind of central clock for managing movements. "
$wsRough.Rows.Item(14).RowHeight = 30
# Minor column width tweaks
$wsRough.Columns.Item(1).ColumnWidth = 68.14285714285714
$wsRough.Columns.Item(2).ColumnWidth = 51.642857142857146

# --- Sheet "Iterations" ---
$wsIter = $wb.Worksheets.Item("Iterations")
$wsIter.Range("B5").Value = "Randomally generate the map"
$wsIter.Columns.Item(1).ColumnWidth = 18.785714285714285

# --- Sheet "Map" ---
$wsMap = $wb.Worksheets.Item("Map")
$wsMap.Activate() | Out-Null
$wsMap.Range("F49:F55").Select() | Out-Null

# --- Make "Rough Guide" the active/selected sheet & cell (activate last so it "wins") ---
$wsRough.Activate() | Out-Null
$wsRough.Range("D12").Select() | Out-Null
